$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Oklo Inc. / OKLO)
$ws.Range("D2").Value = 112.47
$ws.Range("E2").Value = 57.9
$ws.Range("F2").Value = 26.76
$ws.Range("H2").Value = 70
$ws.Range("I2").Value = 73
$ws.Range("J2").Value = 80
$ws.Range("K2").Value = 60.4
$ws.Range("M2").Value = "📈 매수 관찰 구간입니다."
$ws.Range("N2").Value = 54.02451352198364
$ws.Range("O2").Value = "⚪ 중립 구간"

# Row 3 (NuScale Power Corporation / SMR)
$ws.Range("D3").Value = 23.15
$ws.Range("E3").Value = 50
$ws.Range("F3").Value = 21.59
$ws.Range("K3").Value = 53.6
$ws.Range("M3").Value = "⛔ 관망하십시오."
$ws.Range("N3").Value = 54.02451352198364
$ws.Range("O3").Value = "⚪ 중립 구간"
